# Insert two new weekly-update rows right above the current row 476,
# shifting all subsequent rows down by two (517 -> 519).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A476:A477").EntireRow.Insert()

# New row 476 - "Primera" quality entry for the newest reporting date
$ws.Cells.Item(476, 1).Value = 6
$ws.Cells.Item(476, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(476, 3).Value = "Metropolitana"
$ws.Cells.Item(476, 4).Value = 44578
$ws.Cells.Item(476, 5).Value = 13
$ws.Cells.Item(476, 6).Value = 100112017
$ws.Cells.Item(476, 7).Value = "Apio"
$ws.Cells.Item(476, 8).Value = "Americana (o)"
$ws.Cells.Item(476, 9).Value = "Primera"
$ws.Cells.Item(476, 10).Value = 2100
$ws.Cells.Item(476, 11).Value = 6000
$ws.Cells.Item(476, 12).Value = 7000
$ws.Cells.Item(476, 13).Value = 6429
$ws.Cells.Item(476, 14).Value = "`$/docena de matas"
$ws.Cells.Item(476, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(476, 16).Value = 1072
$ws.Cells.Item(476, 17).Value = 6
$ws.Cells.Item(476, 18).Value = "Hortaliza"

# New row 477 - "Segunda" quality entry for the same reporting date
$ws.Cells.Item(477, 1).Value = 6
$ws.Cells.Item(477, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(477, 3).Value = "Metropolitana"
$ws.Cells.Item(477, 4).Value = 44578
$ws.Cells.Item(477, 5).Value = 13
$ws.Cells.Item(477, 6).Value = 100112017
$ws.Cells.Item(477, 7).Value = "Apio"
$ws.Cells.Item(477, 8).Value = "Americana (o)"
$ws.Cells.Item(477, 9).Value = "Segunda"
$ws.Cells.Item(477, 10).Value = 500
$ws.Cells.Item(477, 11).Value = 5000
$ws.Cells.Item(477, 12).Value = 5000
$ws.Cells.Item(477, 13).Value = 5000
$ws.Cells.Item(477, 14).Value = "`$/docena de matas"
$ws.Cells.Item(477, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(477, 16).Value = 833
$ws.Cells.Item(477, 17).Value = 6
$ws.Cells.Item(477, 18).Value = "Hortaliza"
